# StaticData 외래키 테스트 (#73)
# Adds a new "MultiForeignTest" worksheet (multi-level / nested foreign-key
# test fixture) after the last existing sheet ("GroupTest"), fills it with
# the sample data table, and makes it the active sheet/selection.

$wb = $excel.ActiveWorkbook

# Add the new sheet right after the current last sheet so it lands at the
# end of the tab strip (TargetTest, ClassListTest, TypeTest, GroupTest, *MultiForeignTest*).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "MultiForeignTest"

# A1 carries the conventional "top-left data cell" marker used by every
# other sheet in this workbook (TargetTest/ClassListTest -> "B2",
# TypeTest -> "C10", GroupTest -> "C2").
$ws.Range("A1").Value = "D5"

# Header row for the nested-foreign-key sample table.
$ws.Range("D5").Value = "Id"
$ws.Range("E5").Value = "TargetId"
$ws.Range("F5").Value = "Info"

# Data rows.
$ws.Range("D6").Value = 5000
$ws.Range("E6").Value = 1001
$ws.Range("F6").Value = "중첩 외래키"

$ws.Range("D7").Value = 5001
$ws.Range("E7").Value = 1003
$ws.Range("F7").Value = "테스트"

# Make the new sheet active with the same selection Excel recorded.
$ws.Activate() | Out-Null
$ws.Range("A1:G8").Select() | Out-Null
